# Add a new leading "Receipt Id" column to the Expensify export sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:E data right into B:F, opening up a blank column A.
$ws.Columns("A").Insert()

# New header + sequential receipt id values (1..9) for the inserted column.
$ws.Range("A1").Value = "Receipt Id"
for ($i = 2; $i -le 10; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# The hyperlinks that used to live in column E now belong in column F; the
# engine does not re-target hyperlink ranges when a column is inserted, so
# recreate them explicitly (preserving address + destination row).
$hyperlinkAddresses = @{}
foreach ($h in $ws.Hyperlinks) {
    $hyperlinkAddresses[$h.Range.Address()] = $h.Address
}
$ws.Hyperlinks.Delete()
foreach ($cellRef in $hyperlinkAddresses.Keys) {
    $rowNum = $cellRef.Substring(3)
    $target = "F" + $rowNum
    $ws.Hyperlinks.Add($ws.Range($target), $hyperlinkAddresses[$cellRef]) | Out-Null
}
# Re-apply the built-in Hyperlink cell style so the moved cells keep looking
# like the rest of the hyperlink column instead of the engine's fresh style.
foreach ($cellRef in $hyperlinkAddresses.Keys) {
    $rowNum = $cellRef.Substring(3)
    $target = "F" + $rowNum
    $ws.Range($target).Style = "Hyperlink"
}

# Column widths: new "Receipt Id" column and the (slightly) re-sized old
# Timestamp column that is now column B.
$ws.Columns("A").ColumnWidth = 14.8
$ws.Columns("B").ColumnWidth = 17.8

# Selection, matching the saved workbook's last active cell.
$ws.Range("N12").Select()
